# Generate Report for Handoff
#
# Updates the localization-status report for the newly-handed-off files:
#   - 3e97455e-a415-4055-babc-fbb0ad0f56fe.md   (row 8)
#   - 4323503f-1a1c-447b-9dd5-12ec1d8b9db2.md   (row 9)
#   - 6bbd5545-dd48-40fd-b13a-e20c8e3841fa.md   (row 11)
#   - d96fcc9e-cde2-480a-9772-747c00657f31.md   (row 12)
#   - ea51c7ba-0a54-48a0-bc9d-890c6710f457.md   (row 13)
#   - ff4e70b9-c868-452f-9fa7-0872e2008cde.md   (row 14)
#
# For each of those rows, on both locale sheets, the Priority column now
# reports "ht" (handoff type), and the handoff timestamps advance to the
# moment the new xliff packages were generated.

$wb = $excel.ActiveWorkbook

$rows = @(8, 9, 11, 12, 13, 14)

# --- Overview sheet: "Latest HO Xliff Generate Date" (column G) ---
$overview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $overview.Range("G$r").Value = "2016-08-17 22:21:08"
}

# --- zh-cn sheet: Priority (E) + Latest Handoff Datetime (H) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $zhcn.Range("E$r").Value = "ht"
    $zhcn.Range("H$r").Value = "2016-08-17 22:20:58"
}

# --- de-de sheet: Priority (E) + Latest Handoff Datetime (H) ---
$dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $dede.Range("E$r").Value = "ht"
    $dede.Range("H$r").Value = "2016-08-17 22:21:08"
}
